# add: workout page i18n
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -------------------------------------------------------------------------
# Reference cells used purely as format donors (their own formatting is not
# touched by this script, so they reliably keep the original style classes
# already present in the workbook: A2 = plain/no style, A27 = separator
# style, A28 = bold "section label" style).
# -------------------------------------------------------------------------

# --- Row 16: new "error" i18n block (A16 already carries the separator style) ---
$ws.Range("A16").Value = "error"
$ws.Range("B16").Value = "wrong.schedule"
$ws.Range("C16").Value = "Don't have schedule. Please, check again"
$ws.Range("D16").Value = "잘못된 접근입니다. 스케줄을 확인해주세요."

# --- Row 17 used to hold the "alert / modal.title" block -> now blank ---
$ws.Range("B17").Copy()
$ws.Range("A17").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A17:D17").ClearContents()

# --- Row 18 stays blank ---

# --- Row 19 now holds what used to be in row 17 (alert/modal.title) with plain formatting ---
$ws.Range("A2").Copy()
$ws.Range("A19").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A19").Value = "alert"
$ws.Range("B19").Value = "modal.title"
$ws.Range("C19").Value = "Alert"
$ws.Range("D19").Value = "알림"

# --- Row 20 used to hold exerciseDataInfo/instructions -> now blank (plain formatting) ---
$ws.Range("A2").Copy()
$ws.Range("A20").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A20:D20").ClearContents()

# --- Row 21 now holds exerciseDataInfo/history (bold style, previously held by old row 19) ---
$ws.Range("A28").Copy()
$ws.Range("A21").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A21").Value = "exerciseDataInfo"
$ws.Range("B21").Value = "history"
$ws.Range("C21").Value = "History"
$ws.Range("D21").Value = "최근 이력"

# --- Row 22 now holds exerciseDataInfo/instructions (bold style, stays on row 22) ---
$ws.Range("A22").Value = "exerciseDataInfo"
$ws.Range("B22").Value = "instructions"
$ws.Range("C22").Value = "Instructions"
$ws.Range("D22").Value = "설명"

# --- Row 24 now holds scheduleList/schedule.row.title (bold style, previously on row 22) ---
$ws.Range("A28").Copy()
$ws.Range("A24").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A24").Value = "scheduleList"
$ws.Range("B24").Value = "schedule.row.title"
$ws.Range("C24").Value = "Schedule {{n}}"
$ws.Range("D24").Value = "스케줄 {{n}}"

# --- Row 25 becomes a plain separator row (previously bold-styled blank row) ---
$ws.Range("A27").Copy()
$ws.Range("A25").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A25:D25").ClearContents()

# Row 26 (scheduleList/schedule.actionBtn.start) is unchanged.

# --- New workout i18n rows 39-41 ---
$ws.Range("A39").Value = "workout"
$ws.Range("B39").Value = "actionBtn.start"
$ws.Range("C39").Value = "Start"
$ws.Range("D39").Value = "시작하기"

$ws.Range("A40").Value = "workout"
$ws.Range("B40").Value = "actionBtn.pause"
$ws.Range("C40").Value = "Pause"
$ws.Range("D40").Value = "일시정지"

$ws.Range("A41").Value = "workout"
$ws.Range("B41").Value = "actionBtn.finish"
$ws.Range("C41").Value = "Finish"
$ws.Range("D41").Value = "종료"

# --- Column widths: column C widened (to mirror the bestFit width applied in Excel) ---
$ws.Columns("C").ColumnWidth = 34.7

# --- Sheet view: scroll position + active selection ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 30
$ws.Range("A41").Select()
